$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell D5: append a new red-colored run "-Tạo wire Frame cho web"
# after a trailing newline to the existing "Thiết kế giao diện..." text.
$d5 = $ws.Range("D5")
$d5Old = $d5.Text
$d5Suffix = "-Tạo wire Frame cho web"
$d5.Value = $d5Old + "`n" + $d5Suffix

$d5StartPos = $d5Old.Length + 2
$d5Len = $d5Suffix.Length
$d5Chars = $d5.Characters($d5StartPos, $d5Len)
$d5Chars.Font.Color = 255
$d5Chars.Font.Name = "Times New Roman"
$d5Chars.Font.Size = 12

# --- Cell H5: append "- Deploy lên web surge.sh" as a new line to the
# existing "Bàn giao code..." text.
$h5 = $ws.Range("H5")
$h5Old = $h5.Text
$h5.Value = $h5Old + "`n" + "- Deploy lên web surge.sh"

# --- Selection / view: last active cell ends up being J6.
$ws.Range("J6").Select() | Out-Null
